# Fruta / hortaliza, semanal
# The weekly refresh re-shuffles the per-record Fecha/Volumen/Precio/Origen
# values across the existing rows (Mercado/Producto/etc. columns are
# identical for every row on this sheet, so only D, M, N, O, P, R, S move).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Snapshot of the "before" values for the columns that get reshuffled,
# keyed by the row they currently live in.
$data = @{
    2  = @{ D = 44592; M = 5;   N = 7500; O = 7500; P = 7500; R = "Región de La Araucanía"; S = 7500 }
    3  = @{ D = 44175; M = 40;  N = 5000; O = 5000; P = 5000; R = "Provincia de Curicó";    S = 5000 }
    4  = @{ D = 44323; M = 20;  N = 3200; O = 3200; P = 3200; R = "Región de La Araucanía"; S = 3200 }
    5  = @{ D = 44999; M = 25;  N = 2500; O = 2500; P = 2500; R = "Región de La Araucanía"; S = 2500 }
    6  = @{ D = 44215; M = 65;  N = 2800; O = 2800; P = 2800; R = "Región de La Araucanía"; S = 2800 }
    7  = @{ D = 44214; M = 50;  N = 1800; O = 1800; P = 1800; R = "Región de La Araucanía"; S = 1800 }
    8  = @{ D = 44551; M = 120; N = 4500; O = 4500; P = 4500; R = "Región de O'Higgins";    S = 4500 }
    9  = @{ D = 44616; M = 200; N = 3200; O = 3200; P = 3200; R = "Región de La Araucanía"; S = 3200 }
    10 = @{ D = 44176; M = 20;  N = 3000; O = 3000; P = 3000; R = "Región de O'Higgins";    S = 3000 }
    11 = @{ D = 44998; M = 20;  N = 2500; O = 2500; P = 2500; R = "Región de La Araucanía"; S = 2500 }
    12 = @{ D = 44567; M = 80;  N = 2400; O = 2400; P = 2400; R = "Región de La Araucanía"; S = 2400 }
    13 = @{ D = 44574; M = 200; N = 3000; O = 3000; P = 3000; R = "Región de La Araucanía"; S = 3000 }
}

# Maps each destination row to the row whose "before" data it should now
# hold (a permutation of rows 2-13; rows 4 and 12 are fixed points).
$mapping = @{
    2  = 10
    3  = 6
    4  = 4
    5  = 3
    6  = 11
    7  = 5
    8  = 9
    9  = 8
    10 = 13
    11 = 7
    12 = 12
    13 = 2
}

foreach ($destRow in ($mapping.Keys | Sort-Object)) {
    $src = $data[$mapping[$destRow]]

    $ws.Range("D$destRow").Value = $src.D
    $ws.Range("M$destRow").Value = $src.M
    $ws.Range("N$destRow").Value = $src.N
    $ws.Range("O$destRow").Value = $src.O
    $ws.Range("P$destRow").Value = $src.P
    $ws.Range("R$destRow").Value = $src.R
    $ws.Range("S$destRow").Value = $src.S
}
